# Finished Week 13 logging
$wb = $excel.ActiveWorkbook

# Offense sheet ("OFF") - Home row (row 2) updated totals
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 163
$wsOff.Range("C2").Value = 104
$wsOff.Range("D2").Value = 36
$wsOff.Range("E2").Value = 13
$wsOff.Range("F2").Value = 2
$wsOff.Range("G2").Value = 5

# Defense sheet ("DEF") - Home row (row 2) updated totals
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 206
$wsDef.Range("C2").Value = 146
$wsDef.Range("D2").Value = 65
$wsDef.Range("E2").Value = 30
$wsDef.Range("G2").Value = 3
